# Refresh market-price / profit columns (H-N) across all job sheets.
# Values correspond to a scheduled market-data pull; see commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 498
$ws.Range("J9").Value = 233.33333
$ws.Range("L9").Value = 233.33333
$ws.Range("N9").Value = -571.3333299999999
$ws.Range("H17").Value = 466544.3
$ws.Range("J17").Value = 477272.1
$ws.Range("L17").Value = 1431816.3
$ws.Range("N17").Value = -1432152.3
$ws.Range("H18").Value = 336.6
$ws.Range("I18").Value = 336.6
$ws.Range("K18").Value = 336.6
$ws.Range("M18").Value = -52.60000000000002
$ws.Range("H69").Value = 9917.593000000001
$ws.Range("I69").Value = 8000
$ws.Range("J69").Value = 9991.346
$ws.Range("K69").Value = 24000
$ws.Range("L69").Value = 29974.038
$ws.Range("M69").Value = -23126
$ws.Range("N69").Value = -31722.038
$ws.Range("H72").Value = 9917.593000000001
$ws.Range("I72").Value = 8000
$ws.Range("J72").Value = 9991.346
$ws.Range("K72").Value = 72000
$ws.Range("L72").Value = 89922.114
$ws.Range("M72").Value = -67632
$ws.Range("N72").Value = -98658.114
$ws.Range("H74").Value = 6811.727
$ws.Range("I74").Value = 10485.8
$ws.Range("J74").Value = 3750
$ws.Range("K74").Value = 10485.8
$ws.Range("L74").Value = 3750
$ws.Range("M74").Value = -9549.799999999999
$ws.Range("N74").Value = -5622
$ws.Range("H77").Value = 6811.727
$ws.Range("I77").Value = 10485.8
$ws.Range("J77").Value = 3750
$ws.Range("K77").Value = 52429
$ws.Range("L77").Value = 18750
$ws.Range("M77").Value = -47749
$ws.Range("N77").Value = -28110
$ws.Range("H116").Value = 71431656
$ws.Range("I116").Value = 100003000
$ws.Range("K116").Value = 100003000
$ws.Range("M116").Value = -99999558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4695.05
$ws.Range("I45").Value = 4400.0625
$ws.Range("K45").Value = 4400.0625
$ws.Range("M45").Value = -4023.0625
$ws.Range("H132").Value = 1118382.9
$ws.Range("I132").Value = 1672574.4
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 5017723.199999999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -5015193.199999999
$ws.Range("N132").Value = -35060
$ws.Range("H139").Value = 100715
$ws.Range("J139").Value = 100715
$ws.Range("L139").Value = 100715
$ws.Range("N139").Value = -110995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2944.5
$ws.Range("I20").Value = 2871.9375
$ws.Range("J20").Value = 3138
$ws.Range("K20").Value = 2871.9375
$ws.Range("L20").Value = 3138
$ws.Range("M20").Value = -2624.9375
$ws.Range("N20").Value = -3632
$ws.Range("H64").Value = 800.5
$ws.Range("J64").Value = 825.5
$ws.Range("L64").Value = 825.5
$ws.Range("N64").Value = -1275.5
$ws.Range("H67").Value = 800.5
$ws.Range("J67").Value = 825.5
$ws.Range("L67").Value = 825.5
$ws.Range("N67").Value = -2385.5
$ws.Range("H105").Value = 3359.5833
$ws.Range("I105").Value = 3229.4443
$ws.Range("K105").Value = 3229.4443
$ws.Range("M105").Value = -1482.4443
$ws.Range("H134").Value = 3926118.2
$ws.Range("I134").Value = 4448206.5
$ws.Range("J134").Value = 10457
$ws.Range("K134").Value = 13344619.5
$ws.Range("L134").Value = 31371
$ws.Range("M134").Value = -13342084.5
$ws.Range("N134").Value = -36441

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3501.4243
$ws.Range("I31").Value = 1702.85
$ws.Range("J31").Value = 4283.413
$ws.Range("K31").Value = 1702.85
$ws.Range("L31").Value = 4283.413
$ws.Range("M31").Value = -1407.85
$ws.Range("N31").Value = -4873.413
$ws.Range("H34").Value = 3501.4243
$ws.Range("I34").Value = 1702.85
$ws.Range("J34").Value = 4283.413
$ws.Range("K34").Value = 1702.85
$ws.Range("L34").Value = 4283.413
$ws.Range("M34").Value = -1500.85
$ws.Range("N34").Value = -4687.413
$ws.Range("H99").Value = 3460.6924
$ws.Range("I99").Value = 3056.1428
$ws.Range("J99").Value = 3932.6667
$ws.Range("K99").Value = 3056.1428
$ws.Range("L99").Value = 3932.6667
$ws.Range("M99").Value = -1558.1428
$ws.Range("N99").Value = -6928.6667
$ws.Range("H126").Value = 3460.6924
$ws.Range("I126").Value = 3056.1428
$ws.Range("J126").Value = 3932.6667
$ws.Range("K126").Value = 9168.428400000001
$ws.Range("L126").Value = 11798.0001
$ws.Range("M126").Value = -6698.428400000001
$ws.Range("N126").Value = -16738.0001
$ws.Range("H132").Value = 4073.0417
$ws.Range("I132").Value = 4186.294
$ws.Range("K132").Value = 12558.882
$ws.Range("M132").Value = -10028.882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 41.666668
$ws.Range("J2").Value = 113.375
$ws.Range("L2").Value = 680.25
$ws.Range("N2").Value = -906.25
$ws.Range("H4").Value = 130168264
$ws.Range("H34").Value = 2662.389
$ws.Range("J34").Value = 3505.077
$ws.Range("L34").Value = 10515.231
$ws.Range("N34").Value = -10683.231
$ws.Range("H38").Value = 62.23077
$ws.Range("J38").Value = 54.5
$ws.Range("L38").Value = 163.5
$ws.Range("N38").Value = -857.5
$ws.Range("H44").Value = 928.8570999999999
$ws.Range("I44").Value = 1069.6
$ws.Range("J44").Value = 577
$ws.Range("K44").Value = 3208.8
$ws.Range("L44").Value = 1731
$ws.Range("M44").Value = -2810.8
$ws.Range("N44").Value = -2527
$ws.Range("H55").Value = 4119.0835
$ws.Range("I55").Value = 1966.3334
$ws.Range("J55").Value = 4836.6665
$ws.Range("K55").Value = 5899.0002
$ws.Range("L55").Value = 14509.9995
$ws.Range("M55").Value = -5722.0002
$ws.Range("N55").Value = -14863.9995
$ws.Range("H56").Value = 5897.6294
$ws.Range("I56").Value = 5897.6294
$ws.Range("K56").Value = 5897.6294
$ws.Range("M56").Value = -5367.6294
$ws.Range("H68").Value = 2274.6956
$ws.Range("I68").Value = 1851.9
$ws.Range("J68").Value = 2599.923
$ws.Range("K68").Value = 5555.700000000001
$ws.Range("L68").Value = 7799.768999999999
$ws.Range("M68").Value = -4744.700000000001
$ws.Range("N68").Value = -9421.769
$ws.Range("H71").Value = 2274.6956
$ws.Range("I71").Value = 1851.9
$ws.Range("J71").Value = 2599.923
$ws.Range("K71").Value = 16667.1
$ws.Range("L71").Value = 23399.307
$ws.Range("M71").Value = -12611.1
$ws.Range("N71").Value = -31511.307
$ws.Range("H119").Value = 3063.3635
$ws.Range("I119").Value = 2337.125
$ws.Range("K119").Value = 7011.375
$ws.Range("M119").Value = -2173.375
$ws.Range("H132").Value = 3040.7727
$ws.Range("J132").Value = 3233.3333
$ws.Range("L132").Value = 29099.9997
$ws.Range("N132").Value = -34159.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 74950
$ws.Range("J100").Value = 74950
$ws.Range("L100").Value = 74950
$ws.Range("N100").Value = -77114
$ws.Range("H113").Value = 39418.043
$ws.Range("I113").Value = 6535.7144
$ws.Range("K113").Value = 6535.7144
$ws.Range("M113").Value = -4365.7144
$ws.Range("H118").Value = 111999
$ws.Range("J118").Value = 111999
$ws.Range("L118").Value = 111999
$ws.Range("N118").Value = -115313
$ws.Range("H126").Value = 4381.6665
$ws.Range("I126").Value = 4328.3335
$ws.Range("J126").Value = 4488.3335
$ws.Range("K126").Value = 12985.0005
$ws.Range("L126").Value = 13465.0005
$ws.Range("M126").Value = -10515.0005
$ws.Range("N126").Value = -18405.0005
$ws.Range("H135").Value = 206924
$ws.Range("J135").Value = 206924
$ws.Range("L135").Value = 206924
$ws.Range("N135").Value = -217064

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8406.588
$ws.Range("I7").Value = 7037.5454
$ws.Range("J7").Value = 10916.5
$ws.Range("K7").Value = 7037.5454
$ws.Range("L7").Value = 10916.5
$ws.Range("M7").Value = -6925.5454
$ws.Range("N7").Value = -11140.5
$ws.Range("H16").Value = 2505.9
$ws.Range("I16").Value = 2508.875
$ws.Range("K16").Value = 2508.875
$ws.Range("M16").Value = -2338.875
$ws.Range("H68").Value = 4095
$ws.Range("I68").Value = 3122.1428
$ws.Range("K68").Value = 3122.1428
$ws.Range("M68").Value = -2373.1428
$ws.Range("H71").Value = 4095
$ws.Range("I71").Value = 3122.1428
$ws.Range("K71").Value = 15610.714
$ws.Range("M71").Value = -11866.714
$ws.Range("H126").Value = 8406.588
$ws.Range("I126").Value = 7037.5454
$ws.Range("J126").Value = 10916.5
$ws.Range("K126").Value = 21112.6362
$ws.Range("L126").Value = 32749.5
$ws.Range("M126").Value = -18642.6362
$ws.Range("N126").Value = -37689.5
$ws.Range("H134").Value = 99999
$ws.Range("J134").Value = 99999
$ws.Range("L134").Value = 99999
$ws.Range("N134").Value = -110139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 54517
$ws.Range("H73").Value = 54517
